$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D column (price) to remain text for rows whose new value
# would otherwise be auto-detected as a number by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '58.689.78'
$ws.Range("E2").Value = '  -3.86%  '
$ws.Range("D3").Value = '2.609.97'
$ws.Range("E3").Value = '  -2.32%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '509.41'
$ws.Range("E5").Value = '  -3.95%  '
$ws.Range("D6").Value = '146.67'
$ws.Range("E6").Value = '  -5.81%  '
$ws.Range("E8").Value = '  -1.58%  '
$ws.Range("D9").Value = '2.633.22'
$ws.Range("E9").Value = '  -1.71%  '
$ws.Range("D10").Value = '6.39'
$ws.Range("E10").Value = '  -1.50%  '
$ws.Range("D11").Value = '0.105'
$ws.Range("E11").Value = '  -4.00%  '
$ws.Range("D12").Value = '0.338'
$ws.Range("D13").Value = '0.128'
$ws.Range("E13").Value = '  -0.97%  '
$ws.Range("D14").Value = '3.069.21'
$ws.Range("E14").Value = '  -1.70%  '
$ws.Range("D15").Value = '57.965.53'
$ws.Range("E15").Value = '  -5.03%  '
$ws.Range("D16").Value = '21.20'
$ws.Range("E16").Value = '  -4.09%  '
$ws.Range("E17").Value = '  -3.38%  '
$ws.Range("D18").Value = '2.617.84'
$ws.Range("E18").Value = '  -2.04%  '
$ws.Range("D19").Value = '4.59'
$ws.Range("E19").Value = '  -4.06%  '
$ws.Range("D20").Value = '345.50'
$ws.Range("E20").Value = '  -2.61%  '
$ws.Range("D21").Value = '10.42'
$ws.Range("E21").Value = '  -2.53%  '
$ws.Range("D22").Value = '6.18'
$ws.Range("E22").Value = '  -2.50%  '
$ws.Range("E23").Value = '  -0.15%  '
$ws.Range("E24").Value = '  -1.50%  '
$ws.Range("E25").Value = '  -1.98%  '
$ws.Range("D26").Value = '2.720.46'
$ws.Range("E26").Value = '  -2.20%  '
$ws.Range("E27").Value = '  -0.31%  '
$ws.Range("D28").Value = '0.161'
$ws.Range("E28").Value = '  -4.27%  '
$ws.Range("D29").Value = '0.0₃0822'
$ws.Range("E29").Value = '  -4.11%  '
$ws.Range("D30").Value = '7.06'
$ws.Range("E30").Value = '  -3.73%  '
$ws.Range("E31").Value = '  -0.14%  '
$ws.Range("D32").Value = '6.24'
$ws.Range("E32").Value = '  +0.95%  '
$ws.Range("D33").Value = '19.00'
$ws.Range("E33").Value = '  -2.73%  '
$ws.Range("E34").Value = '  -4.53%  '
$ws.Range("D35").Value = '149.29'
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("D36").Value = '0.989'
$ws.Range("E36").Value = '  +11.14%  '
$ws.Range("D37").Value = '4.02'
$ws.Range("E37").Value = '  -2.57%  '
$ws.Range("E38").Value = '  -4.32%  '
$ws.Range("D39").Value = '0.869'
$ws.Range("E39").Value = '  -5.61%  '
$ws.Range("E40").Value = '  -2.07%  '
$ws.Range("E41").Value = '  -4.41%  '
$ws.Range("D42").Value = '3.64'
$ws.Range("E42").Value = '  -3.78%  '
$ws.Range("D43").Value = '289.62'
$ws.Range("E43").Value = '  -5.34%  '
$ws.Range("E44").Value = '  -4.25%  '
$ws.Range("E45").Value = '  -2.07%  '
$ws.Range("D46").Value = '0.994'
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("D47").Value = '19.70'
$ws.Range("E47").Value = '  -3.47%  '
$ws.Range("D48").Value = '0.0540'
$ws.Range("E48").Value = '  -4.35%  '
$ws.Range("D49").Value = '4.75'
$ws.Range("E49").Value = '  -3.93%  '
$ws.Range("D50").Value = '0.0231'
$ws.Range("E50").Value = '  -4.70%  '
$ws.Range("E51").Value = '  -0.93%  '

# Restore the original (default/"Normal") style on column D so no
# unintended cell-level formatting change is left behind.
$ws.Range("D2:D51").Style = "Normal"
